$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s) Responsável(eis)*") {
        $target = $p
        break
    }
}

# Insert a new paragraph right after it.
$newRange = $target.Range.InsertParagraphAfter()

# The newly created paragraph is the one following the target paragraph.
$newPara = $target.Next()
$newPara.Range.Text = "7455355 - Robson da Silva Rocha"
$newPara.Style = "ListBullet"
